# Applies: "added date and bottom row border"
#  1. Row 2 becomes the report "date range" line -- non-bold, smaller
#     centered text showing &=display.ReportDateRange, row height reverts
#     to the sheet default (no more custom height).
#  2. The header/formula row (row 4, A4:K4) gets a thin gray bottom
#     border added beneath the existing borders, with a couple of small
#     companion cells (M4, K5, L5) carrying matching border fragments.
#  3. Selection moves to F6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a single border edge with an explicit color so the
# serialized XML gets a concrete <color rgb="..."/> (COM doesn't persist
# ThemeColor/TintAndShade on Border objects in this host).
function Set-Edge($range, $edgeIndex, $colorBgr) {
    $edge = $range.Borders.Item($edgeIndex)
    $edge.LineStyle = 1   # xlContinuous
    $edge.Weight = 2      # xlThin
    $edge.Color = $colorBgr
}

$xlEdgeLeft = 7
$xlEdgeTop = 8
$xlEdgeBottom = 9
$xlEdgeRight = 10

$black = 0
$grayThemeish = 13553360   # BGR for D0CECE == theme(2) tint(-0.099978637043366805)

# --- 1. Date-range line (row 2) -----------------------------------------
$ws.Range("A2").Value = "&=display.ReportDateRange"

$row2Range = $ws.Range("A2:L2")
$row2Range.Font.Bold = $false
$row2Range.Font.Size = 9
$row2Range.HorizontalAlignment = -4108   # xlCenter

# Drop the custom row height the title row used to need -- AutoFit lets
# Excel fall back to the sheet default (12.5) with no customHeight flag.
$ws.Rows.Item(2).AutoFit()

# --- 2. Bottom border under the formula/header row (row 4) -------------
Set-Edge $ws.Range("A4:K4") $xlEdgeBottom $grayThemeish

# Small filler cell right after L4 picking up the left divider border
Set-Edge $ws.Range("M4") $xlEdgeLeft $grayThemeish
$ws.Range("M4").VerticalAlignment = -4160   # xlTop

# New row 5 continues the divider under K4:L4 with a thin top border
Set-Edge $ws.Range("K5:L5") $xlEdgeTop $grayThemeish

# --- 3. Cosmetic: move the active selection ------------------------------
$ws.Range("F6").Select()
